# Apply updated crypto price/volume values (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.768.05'
$ws.Range("E2").Value = '  -4.63%  '

$ws.Range("D3").Value = '2.322.50'
$ws.Range("E3").Value = '  -6.19%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.86'
$ws.Range("E5").Value = '  -4.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '83.58'
$ws.Range("E6").Value = '  -9.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.526'
$ws.Range("E7").Value = '  -4.58%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("E9").Value = '  -6.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0812'
$ws.Range("E10").Value = '  -5.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.70'
$ws.Range("E11").Value = '  -10.27%  '

$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = '2.693.55'
$ws.Range("E13").Value = '  -5.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.35'
$ws.Range("E14").Value = '  -7.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.68'
$ws.Range("E15").Value = '  -5.37%  '

$ws.Range("D16").Value = '2.329.80'
$ws.Range("E16").Value = '  -5.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.742'
$ws.Range("E17").Value = '  -6.67%  '

$ws.Range("D18").Value = '39.802.22'
$ws.Range("E18").Value = '  -4.37%  '

$ws.Range("D19").Value = '0.0₃0894'
$ws.Range("E19").Value = '  -5.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.02'
$ws.Range("E20").Value = '  -6.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.42'
$ws.Range("E21").Value = '  -4.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.45'
$ws.Range("E22").Value = '  -7.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.92'
$ws.Range("E23").Value = '  -2.83%  '

$ws.Range("E24").Value = '  -8.30%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.78'
$ws.Range("E26").Value = '  -8.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.20'
$ws.Range("E27").Value = '  -7.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -2.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("E29").Value = '  -6.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.85'
$ws.Range("E30").Value = '  -7.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '150.03'
$ws.Range("E31").Value = '  -4.76%  '

$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.07'
$ws.Range("E33").Value = '  -6.81%  '

$ws.Range("E36").Value = '  -3.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.73'
$ws.Range("E37").Value = '  -5.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0985'
$ws.Range("E38").Value = '  -5.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.53'
$ws.Range("E39").Value = '  -9.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.67'
$ws.Range("E40").Value = '  -9.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.75'
$ws.Range("E41").Value = '  -6.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.32'
$ws.Range("E42").Value = '  -4.76%  '

$ws.Range("D43").Value = '1.938.65'
$ws.Range("E43").Value = '  -2.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0263'
$ws.Range("E44").Value = '  -7.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.33'
$ws.Range("E45").Value = '  -7.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.24'
$ws.Range("E46").Value = '  -2.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.62'
$ws.Range("E47").Value = '  -12.12%  '

$ws.Range("D48").Value = '2.569.14'
$ws.Range("E48").Value = '  -6.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '91.37'
$ws.Range("E49").Value = '  -6.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.28'
$ws.Range("E50").Value = '  -7.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.74'
$ws.Range("E51").Value = '  -6.99%  '

# Row 34/35: Hedera and WEMIXToken swapped positions with updated price/volume
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.41'
$ws.Range("E34").Value = '  -5.71%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0715'
$ws.Range("E35").Value = '  -6.48%  '

